$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02506566666666667
$ws.Range("H2").Value = 0.075197
$ws.Range("M2").Value = 12.40685866666667
$ws.Range("N2").Value = 37.220576
$ws.Range("O2").Value = 0.1720325859617629
$ws.Range("P2").Value = 0.1720325859617629
$ws.Range("Q2").Value = 0.3109861837191111
$ws.Range("R2").Value = 2.798875653472
$ws.Range("S2").Value = 0.1720325859617629
$ws.Range("T2").Value = 0.1720325859617629

# Row 3
$ws.Range("G3").Value = 0.02506566666666667
$ws.Range("H3").Value = 0.075197
$ws.Range("O3").Value = 0.6097142007069145
$ws.Range("P3").Value = 0.6097142007069145
$ws.Range("Q3").Value = 1.102190561033222
$ws.Range("R3").Value = 9.919715049298999
$ws.Range("S3").Value = 0.6097142007069145
$ws.Range("T3").Value = 0.6097142007069145

# Row 4
$ws.Range("G4").Value = 0.02506566666666667
$ws.Range("H4").Value = 0.075197
$ws.Range("N4").Value = 47.220765
$ws.Range("O4").Value = 0.2182532133313226
$ws.Range("P4").Value = 0.2182532133313226
$ws.Range("Q4").Value = 0.3945399850783333
$ws.Range("R4").Value = 3.550859865705
$ws.Range("S4").Value = 0.2182532133313226
$ws.Range("T4").Value = 0.2182532133313226
